$wb = $excel.ActiveWorkbook

$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIR.Range("A225:F225").NumberFormat = "@"
$wsPIR.Range("A225").Value = "2026-01-30"
$wsPIR.Range("B225").Value = "17:25:41"
$wsPIR.Range("C225").Value = "17:00"
$wsPIR.Range("D225").Value = "Bathroom"
$wsPIR.Range("E225").Value = "No Motion"
$wsPIR.Range("F225").Value = "Inactive"

$wsPIR.Range("A226:F226").NumberFormat = "@"
$wsPIR.Range("A226").Value = "2026-01-30"
$wsPIR.Range("B226").Value = "17:25:43"
$wsPIR.Range("C226").Value = "17:00"
$wsPIR.Range("D226").Value = "Bathroom"
$wsPIR.Range("E226").Value = "No Motion"
$wsPIR.Range("F226").Value = "Inactive"

$wsPIR.Range("A227:F227").NumberFormat = "@"
$wsPIR.Range("A227").Value = "2026-01-30"
$wsPIR.Range("B227").Value = "17:25:48"
$wsPIR.Range("C227").Value = "17:00"
$wsPIR.Range("D227").Value = "Bathroom"
$wsPIR.Range("E227").Value = "No Motion"
$wsPIR.Range("F227").Value = "Inactive"

$wsPIR.Range("A228:F228").NumberFormat = "@"
$wsPIR.Range("A228").Value = "2026-01-30"
$wsPIR.Range("B228").Value = "17:25:53"
$wsPIR.Range("C228").Value = "17:00"
$wsPIR.Range("D228").Value = "Bathroom"
$wsPIR.Range("E228").Value = "No Motion"
$wsPIR.Range("F228").Value = "Inactive"

$wsPIR.Range("A229:F229").NumberFormat = "@"
$wsPIR.Range("A229").Value = "2026-01-30"
$wsPIR.Range("B229").Value = "17:25:58"
$wsPIR.Range("C229").Value = "17:00"
$wsPIR.Range("D229").Value = "Bathroom"
$wsPIR.Range("E229").Value = "No Motion"
$wsPIR.Range("F229").Value = "Inactive"

$wsPIR.Range("A230:F230").NumberFormat = "@"
$wsPIR.Range("A230").Value = "2026-01-30"
$wsPIR.Range("B230").Value = "17:26:03"
$wsPIR.Range("C230").Value = "17:00"
$wsPIR.Range("D230").Value = "Bathroom"
$wsPIR.Range("E230").Value = "No Motion"
$wsPIR.Range("F230").Value = "Inactive"

$wsPIR.Range("A231:F231").NumberFormat = "@"
$wsPIR.Range("A231").Value = "2026-01-30"
$wsPIR.Range("B231").Value = "17:26:08"
$wsPIR.Range("C231").Value = "17:00"
$wsPIR.Range("D231").Value = "Bathroom"
$wsPIR.Range("E231").Value = "No Motion"
$wsPIR.Range("F231").Value = "Inactive"

$wsPIR.Range("A232:F232").NumberFormat = "@"
$wsPIR.Range("A232").Value = "2026-01-30"
$wsPIR.Range("B232").Value = "17:26:13"
$wsPIR.Range("C232").Value = "17:00"
$wsPIR.Range("D232").Value = "Bathroom"
$wsPIR.Range("E232").Value = "No Motion"
$wsPIR.Range("F232").Value = "Inactive"

$wsPIR.Range("A233:F233").NumberFormat = "@"
$wsPIR.Range("A233").Value = "2026-01-30"
$wsPIR.Range("B233").Value = "17:26:18"
$wsPIR.Range("C233").Value = "17:00"
$wsPIR.Range("D233").Value = "Bathroom"
$wsPIR.Range("E233").Value = "No Motion"
$wsPIR.Range("F233").Value = "Inactive"

$wsPIR.Range("A234:F234").NumberFormat = "@"
$wsPIR.Range("A234").Value = "2026-01-30"
$wsPIR.Range("B234").Value = "17:26:23"
$wsPIR.Range("C234").Value = "17:00"
$wsPIR.Range("D234").Value = "Bathroom"
$wsPIR.Range("E234").Value = "No Motion"
$wsPIR.Range("F234").Value = "Inactive"

$wsPIR.Range("A235:F235").NumberFormat = "@"
$wsPIR.Range("A235").Value = "2026-01-30"
$wsPIR.Range("B235").Value = "17:26:28"
$wsPIR.Range("C235").Value = "17:00"
$wsPIR.Range("D235").Value = "Bathroom"
$wsPIR.Range("E235").Value = "No Motion"
$wsPIR.Range("F235").Value = "Inactive"

$wsPIR.Range("A236:F236").NumberFormat = "@"
$wsPIR.Range("A236").Value = "2026-01-30"
$wsPIR.Range("B236").Value = "17:26:33"
$wsPIR.Range("C236").Value = "17:00"
$wsPIR.Range("D236").Value = "Bathroom"
$wsPIR.Range("E236").Value = "No Motion"
$wsPIR.Range("F236").Value = "Inactive"

$wsPIR.Range("A237:F237").NumberFormat = "@"
$wsPIR.Range("A237").Value = "2026-01-30"
$wsPIR.Range("B237").Value = "17:26:38"
$wsPIR.Range("C237").Value = "17:00"
$wsPIR.Range("D237").Value = "Bathroom"
$wsPIR.Range("E237").Value = "No Motion"
$wsPIR.Range("F237").Value = "Inactive"

$wsHum = $wb.Worksheets.Item("Humidity")
$wsHum.Range("A149:F149").NumberFormat = "@"
$wsHum.Range("A149").Value = "2026-01-30"
$wsHum.Range("B149").Value = "17:25:42"
$wsHum.Range("C149").Value = "17:00"
$wsHum.Range("D149").Value = "Bathroom"
$wsHum.Range("E149").Value = "86.4%"
$wsHum.Range("F149").Value = "Active"

$wsHum.Range("A150:F150").NumberFormat = "@"
$wsHum.Range("A150").Value = "2026-01-30"
$wsHum.Range("B150").Value = "17:25:53"
$wsHum.Range("C150").Value = "17:00"
$wsHum.Range("D150").Value = "Bathroom"
$wsHum.Range("E150").Value = "87.3%"
$wsHum.Range("F150").Value = "Active"

$wsHum.Range("A151:F151").NumberFormat = "@"
$wsHum.Range("A151").Value = "2026-01-30"
$wsHum.Range("B151").Value = "17:26:08"
$wsHum.Range("C151").Value = "17:00"
$wsHum.Range("D151").Value = "Bathroom"
$wsHum.Range("E151").Value = "86.1%"
$wsHum.Range("F151").Value = "Active"

$wsHum.Range("A152:F152").NumberFormat = "@"
$wsHum.Range("A152").Value = "2026-01-30"
$wsHum.Range("B152").Value = "17:26:13"
$wsHum.Range("C152").Value = "17:00"
$wsHum.Range("D152").Value = "Bathroom"
$wsHum.Range("E152").Value = "87.4%"
$wsHum.Range("F152").Value = "Active"

$wsHum.Range("A153:F153").NumberFormat = "@"
$wsHum.Range("A153").Value = "2026-01-30"
$wsHum.Range("B153").Value = "17:26:18"
$wsHum.Range("C153").Value = "17:00"
$wsHum.Range("D153").Value = "Bathroom"
$wsHum.Range("E153").Value = "87.4%"
$wsHum.Range("F153").Value = "Active"

$wsHum.Range("A154:F154").NumberFormat = "@"
$wsHum.Range("A154").Value = "2026-01-30"
$wsHum.Range("B154").Value = "17:26:28"
$wsHum.Range("C154").Value = "17:00"
$wsHum.Range("D154").Value = "Bathroom"
$wsHum.Range("E154").Value = "87.4%"
$wsHum.Range("F154").Value = "Active"

$wsHum.Range("A155:F155").NumberFormat = "@"
$wsHum.Range("A155").Value = "2026-01-30"
$wsHum.Range("B155").Value = "17:26:33"
$wsHum.Range("C155").Value = "17:00"
$wsHum.Range("D155").Value = "Bathroom"
$wsHum.Range("E155").Value = "87.3%"
$wsHum.Range("F155").Value = "Active"

$wsHum.Range("A156:F156").NumberFormat = "@"
$wsHum.Range("A156").Value = "2026-01-30"
$wsHum.Range("B156").Value = "17:26:38"
$wsHum.Range("C156").Value = "17:00"
$wsHum.Range("D156").Value = "Bathroom"
$wsHum.Range("E156").Value = "87.4%"
$wsHum.Range("F156").Value = "Active"
